# Apply "new sensitivity and calculus" update to the daily model results workbook.
$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------------
# Sheet: "Model Accuracy (-0.1, 0.1, 0.1)"
# Add new columns C:G (Market threshold, Market min, Market max, Recall,
# Precision) and refresh the Accuracy (%) values in column B.
# -----------------------------------------------------------------------
$wsAcc = $wb.Worksheets.Item("Model Accuracy (-0.1, 0.1, 0.1)")

# New header cells, copying the formatting already used by B1 (bold,
# bordered, centered header style) onto the newly added header cells.
$wsAcc.Range("C1").Value = "Market threshold"
$wsAcc.Range("D1").Value = "Market min"
$wsAcc.Range("E1").Value = "Market max"
$wsAcc.Range("F1").Value = "Recall"
$wsAcc.Range("G1").Value = "Precision"

$wsAcc.Range("B1").Copy()
$wsAcc.Range("C1:G1").PasteSpecial(-4122)

# Row 2 - TOTALENERGIES SE
$wsAcc.Range("B2").Value = 32.94621026894865
$wsAcc.Range("C2").Value = 0.05450546436368681
$wsAcc.Range("D2").Value = -15.55441
$wsAcc.Range("E2").Value = 15.06418
$wsAcc.Range("F2").Value = 11.11111111111111
$wsAcc.Range("G2").Value = 1.098901098901099

# Row 3 - FMC CORP
$wsAcc.Range("B3").Value = 27.13936430317848
$wsAcc.Range("C3").Value = 0.009583939973006913
$wsAcc.Range("D3").Value = -19.35264
$wsAcc.Range("E3").Value = 13.70093
$wsAcc.Range("F3").Value = 14.20911528150134
$wsAcc.Range("G3").Value = 25

# Row 4 - BP PLC
$wsAcc.Range("B4").Value = 51.83374083129584
$wsAcc.Range("C4").Value = 0.04158117063764853
$wsAcc.Range("D4").Value = -18.75314
$wsAcc.Range("E4").Value = 23.33066
$wsAcc.Range("F4").Value = 0
$wsAcc.Range("G4").Value = 0

# Row 5 - STORA ENSO
$wsAcc.Range("B5").Value = 41.80929095354523
$wsAcc.Range("C5").Value = 0.02983403801513819
$wsAcc.Range("D5").Value = -12.78028
$wsAcc.Range("E5").Value = 12.42348
$wsAcc.Range("F5").Value = 14.54545454545454
$wsAcc.Range("G5").Value = 6.083650190114068

# Row 6 - BHP GROUP
$wsAcc.Range("B6").Value = 43.94865525672372
$wsAcc.Range("C6").Value = 0.08368817696170747
$wsAcc.Range("D6").Value = -16.47904
$wsAcc.Range("E6").Value = 14.94325
$wsAcc.Range("F6").Value = 50
$wsAcc.Range("G6").Value = 0.3603603603603603

# -----------------------------------------------------------------------
# Sheet: "Confusion Matrix TOTALENERGIES SE (-0.1, 0.1, 0.1)"
# -----------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("Confusion Matrix TOTALENERGIES SE (-0.1, 0.1, 0.1)")
$wsTotal.Range("C3").Value = 531
$wsTotal.Range("D3").Value = 3
$wsTotal.Range("B4").Value = 3
$wsTotal.Range("C4").Value = 478
$wsTotal.Range("D4").Value = 7

# -----------------------------------------------------------------------
# Sheet: "Confusion Matrix FMC CORP (-0.1, 0.1, 0.1)"
# -----------------------------------------------------------------------
$wsFmc = $wb.Worksheets.Item("Confusion Matrix FMC CORP (-0.1, 0.1, 0.1)")
$wsFmc.Range("B2").Value = 53
$wsFmc.Range("C2").Value = 110
$wsFmc.Range("D2").Value = 49
$wsFmc.Range("B3").Value = 124
$wsFmc.Range("C3").Value = 223
$wsFmc.Range("D3").Value = 139
$wsFmc.Range("B4").Value = 196
$wsFmc.Range("C4").Value = 319
$wsFmc.Range("D4").Value = 168

# -----------------------------------------------------------------------
# Sheet: "Confusion Matrix BP PLC (-0.1, 0.1, 0.1)"
# -----------------------------------------------------------------------
$wsBp = $wb.Worksheets.Item("Confusion Matrix BP PLC (-0.1, 0.1, 0.1)")
$wsBp.Range("B3").Value = 21
$wsBp.Range("C3").Value = 830
$wsBp.Range("D3").Value = 24
$wsBp.Range("B4").Value = 19
$wsBp.Range("C4").Value = 619
$wsBp.Range("D4").Value = 18

# -----------------------------------------------------------------------
# Sheet: "Confusion Matrix STORA ENSO (-0.1, 0.1, 0.1)"
# -----------------------------------------------------------------------
$wsStora = $wb.Worksheets.Item("Confusion Matrix STORA ENSO (-0.1, 0.1, 0.1)")
$wsStora.Range("B2").Value = 16
$wsStora.Range("C2").Value = 232
$wsStora.Range("D2").Value = 15
$wsStora.Range("B3").Value = 53
$wsStora.Range("C3").Value = 621
$wsStora.Range("D3").Value = 45
$wsStora.Range("B4").Value = 41
$wsStora.Range("C4").Value = 505
$wsStora.Range("D4").Value = 47

# -----------------------------------------------------------------------
# Sheet: "Confusion Matrix BHP GROUP (-0.1, 0.1, 0.1)"
# -----------------------------------------------------------------------
$wsBhp = $wb.Worksheets.Item("Confusion Matrix BHP GROUP (-0.1, 0.1, 0.1)")
$wsBhp.Range("B2").Value = 2
$wsBhp.Range("C2").Value = 552
$wsBhp.Range("D2").Value = 1
$wsBhp.Range("B3").Value = 1
$wsBhp.Range("C3").Value = 716
$wsBhp.Range("D3").Value = 1
$wsBhp.Range("B4").Value = 1
$wsBhp.Range("C4").Value = 305
$wsBhp.Range("D4").Value = 1
